$wb = $excel.ActiveWorkbook

# --- TestData sheet: remove the DepartureDate / ReturnDate columns (C, D) ---
$testData = $wb.Worksheets.Item("TestData")
$testData.Range("C1:D2").EntireColumn.Delete()
$testData.Range("C1:E1048576").Select()

# --- ObjectRepository sheet: update locators for the new flow ---
$objRepo = $wb.Worksheets.Item("ObjectRepository")

# Rename "Home.ArrivalDate.Xpath" -> "Home.ReturnDate.Xpath"
$objRepo.Range("A10").Value = "Home.ReturnDate.Xpath"

# Update the date-picker locators to the new static values
$objRepo.Range("B9").Value = "(//a[@data-string='2882019'])[1]"
$objRepo.Range("B10").Value = "(//a[@data-string='2792019'])[1]"

# Simplify the leg-OB / leg-IB locators (drop the /text() hops)
$objRepo.Range("B14").Value = '//*[@id="leg-OB"]'
$objRepo.Range("B15").Value = '//*[@id="leg-IB"]'
$objRepo.Range("B16").Value = '//*[@id="leg-OB"]/following::div[1]'
$objRepo.Range("B17").Value = '(//*[@id="leg-IB"]/following::div)[1]'

$objRepo.Range("B15").Select()
